$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Add the new "ParentSoils" worksheet right after "Gullies" (becomes the 3rd/last sheet)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$newSheet.Name = "ParentSoils"

# Header row
$newSheet.Range("A1").Value = "PlanningUnit"
$newSheet.Range("B1").Value = "SoilSource"
$newSheet.Range("C1").Value = "TotalNitrogen"
$newSheet.Range("D1").Value = "TotalCarbon"
$newSheet.Range("E1").Value = "DeltaCarbon"

# Data rows
$newSheet.Range("A2").Value = 17
$newSheet.Range("B2").Value = "Gully"
$newSheet.Range("C2").Value = 0.0775
$newSheet.Range("D2").Value = 2.5025
$newSheet.Range("E2").Formula = "=D2*0.6"

$newSheet.Range("A3").Value = 17
$newSheet.Range("B3").Value = "Hillslope"
$newSheet.Range("C3").Value = 0.13
$newSheet.Range("D3").Value = 1.846666667
$newSheet.Range("E3").Value = 0.443237166

$newSheet.Range("A4").Value = 17
$newSheet.Range("B4").Value = "Riparian"
$newSheet.Range("C4").Value = 0.083992857
$newSheet.Range("D4").Value = 0.98380449
$newSheet.Range("E4").Formula = "=D4*0.6"

$newSheet.Range("A5").Value = 18
$newSheet.Range("B5").Value = "Gully"
$newSheet.Range("C5").Value = 0.104
$newSheet.Range("D5").Value = 1.259333333
$newSheet.Range("E5").Formula = "=D5*0.6"

$newSheet.Range("A6").Value = 18
$newSheet.Range("B6").Value = "Hillslope"
$newSheet.Range("C6").Value = 0.140555556
$newSheet.Range("D6").Value = 1.795212732
$newSheet.Range("E6").Value = 0.41816442

$newSheet.Range("A7").Value = 18
$newSheet.Range("B7").Value = "Riparian"
$newSheet.Range("C7").Value = 0.084071434
$newSheet.Range("D7").Value = 0.985537548
$newSheet.Range("E7").Formula = "=D7*0.6"

$newSheet.Range("A8").Value = 19
$newSheet.Range("B8").Value = "Hillslope"
$newSheet.Range("C8").Value = 0.154
$newSheet.Range("D8").Value = 2.04
$newSheet.Range("E8").Value = 0.381003102

$newSheet.Range("A9").Value = 19
$newSheet.Range("B9").Value = "Riparian"
$newSheet.Range("C9").Value = 0.090518973
$newSheet.Range("D9").Value = 1.229454613
$newSheet.Range("E9").Formula = "=D9*0.6"

$newSheet.Range("A10").Value = 20
$newSheet.Range("B10").Value = "Hillslope"
$newSheet.Range("C10").Value = 0.128428571
$newSheet.Range("D10").Value = 1.929857143
$newSheet.Range("E10").Value = 0.44999999999999996

$newSheet.Range("A11").Value = 20
$newSheet.Range("B11").Value = "Riparian"
$newSheet.Range("C11").Value = 0.104666667
$newSheet.Range("D11").Value = 1.263333333
$newSheet.Range("E11").Formula = "=D11*0.6"

$newSheet.Range("A12").Value = 21
$newSheet.Range("B12").Value = "Hillslope"
$newSheet.Range("C12").Value = 0.146764706
$newSheet.Range("D12").Value = 1.987508403
$newSheet.Range("E12").Value = 0.44999999999999996

$newSheet.Range("A13").Value = 21
$newSheet.Range("B13").Value = "Riparian"
$newSheet.Range("C13").Value = 0.100666667
$newSheet.Range("D13").Value = 1.222488889
$newSheet.Range("E13").Formula = "=D13*0.6"

$newSheet.Range("A14").Value = 22
$newSheet.Range("B14").Value = "Hillslope"
$newSheet.Range("C14").Value = 0.13919797
$newSheet.Range("D14").Value = 1.802498731
$newSheet.Range("E14").Value = 0.44999999999999996

$newSheet.Range("A15").Value = 22
$newSheet.Range("B15").Value = "Riparian"
$newSheet.Range("C15").Value = 0.115166667
$newSheet.Range("D15").Value = 1.386666667
$newSheet.Range("E15").Formula = "=D15*0.6"

$newSheet.Range("A16").Value = 23
$newSheet.Range("B16").Value = "Hillslope"
$newSheet.Range("C16").Value = 0.09
$newSheet.Range("D16").Value = 1.5
$newSheet.Range("E16").Value = 0.444230772

$newSheet.Range("A17").Value = 23
$newSheet.Range("B17").Value = "Riparian"
$newSheet.Range("C17").Value = 0.096666667
$newSheet.Range("D17").Value = 1.177777778
$newSheet.Range("E17").Formula = "=D17*0.6"

# Column widths (approximate best-fit autofit of the final data)
$newSheet.Range("A1:E17").EntireColumn.AutoFit() | Out-Null

# Update selection / active sheet state:
# PlanningUnits loses its tab selection, selection moves to I31
$ws1.Activate()
$ws1.Range("I31").Select()

# ParentSoils becomes the active tab, with selection at G4
$newSheet.Activate()
$newSheet.Range("G4").Select()

Write-Host "done"
